# Update "want to go" counts (column F) across the four sheets to match
# the freshly scraped snapshot referenced in the commit message
# ("Update gh-pages to output generated at 456a3b4").
#
# Sheet order (per xl/workbook.xml): 1=展览, 2=演出, 3=本地生活, 4=全部类型

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F4").Value  = 3214
$ws1.Range("F15").Value = 1068
$ws1.Range("F16").Value = 1003
$ws1.Range("F19").Value = 305
$ws1.Range("F20").Value = 5921
$ws1.Range("F21").Value = 2335
$ws1.Range("F22").Value = 4096
$ws1.Range("F23").Value = 2280
$ws1.Range("F25").Value = 62
$ws1.Range("F26").Value = 62
$ws1.Range("F31").Value = 70
$ws1.Range("F36").Value = 564
$ws1.Range("F41").Value = 286

# --- Sheet 2: 演出 -----------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F6").Value  = 105
$ws2.Range("F9").Value  = 584
$ws2.Range("F18").Value = 128
$ws2.Range("F25").Value = 3810
$ws2.Range("F26").Value = 3810

# --- Sheet 3: 本地生活 ---------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F10").Value = 1215

# --- Sheet 4: 全部类型 ---------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value  = 3214
$ws4.Range("F8").Value  = 1215
$ws4.Range("F12").Value = 105
$ws4.Range("F16").Value = 584
$ws4.Range("F17").Value = 584
$ws4.Range("F20").Value = 1068
$ws4.Range("F26").Value = 305
$ws4.Range("F27").Value = 5921
$ws4.Range("F28").Value = 2335
$ws4.Range("F29").Value = 4096
$ws4.Range("F30").Value = 2280
$ws4.Range("F31").Value = 62
$ws4.Range("F32").Value = 62
$ws4.Range("F35").Value = 70
$ws4.Range("F38").Value = 128
$ws4.Range("F40").Value = 564
$ws4.Range("F46").Value = 287
$ws4.Range("F48").Value = 3810
